# "Generate Report for Archive"
#
# The localization-status report is regenerated: every cell whose status
# was "Ready for handoff" moves to "In Translation", and the two status
# columns on the Overview sheet (and the Status column on each per-locale
# sheet) are re-autofit to the new, shorter text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status columns E (zh-cn) and F (de-de) -----------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

# Columns E and F were autosized to the old "Ready for handoff" text;
# re-fit them to the new, shorter "In Translation" text.
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5

# --- Per-locale sheets: Status column C --------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 12.5
